# Constrain level-1 fields in lightsheet
#
# The "lightsheet" metadata template only needs a single allowed value for
# each of assay_category / assay_type / analyte_class, so trim each lookup
# sheet down to its one surviving entry and point the corresponding data
# validations (columns I, J, K on "Export as TSV") at the now-single-cell
# lookup ranges, updating their error text to match.

$wb = $excel.ActiveWorkbook

# --- assay_category list: keep only "imaging" (row 1); drop the rest ----
$wsCat = $wb.Worksheets.Item("assay_category list")
$wsCat.Range("A2:A4").EntireRow.Delete() | Out-Null

# --- assay_type list: keep a single row, replaced with "Light Sheet" ----
$wsType = $wb.Worksheets.Item("assay_type list")
$wsType.Range("A1").Value = "Light Sheet"
$wsType.Range("A2:A30").EntireRow.Delete() | Out-Null

# --- analyte_class list: keep only "protein" (row 1); drop the rest -----
$wsAnalyte = $wb.Worksheets.Item("analyte_class list")
$wsAnalyte.Range("A1").Value = "protein"
$wsAnalyte.Range("A2:A7").EntireRow.Delete() | Out-Null

# --- Update the data validations on the main sheet ----------------------
$wsMain = $wb.Worksheets.Item("Export as TSV")

$rngI = $wsMain.Range("I2:I1048576")
$rngI.Validation.Modify(3, 1, 1, "'assay_category list'!`$A`$1:`$A`$1")
$rngI.Validation.ErrorTitle = "Value must come from list"
$rngI.Validation.ErrorMessage = "Value must be one of: imaging."

$rngJ = $wsMain.Range("J2:J1048576")
$rngJ.Validation.Modify(3, 1, 1, "'assay_type list'!`$A`$1:`$A`$1")
$rngJ.Validation.ErrorTitle = "Value must come from list"
$rngJ.Validation.ErrorMessage = "Value must be one of: Light Sheet."

$rngK = $wsMain.Range("K2:K1048576")
$rngK.Validation.Modify(3, 1, 1, "'analyte_class list'!`$A`$1:`$A`$1")
$rngK.Validation.ErrorTitle = "Value must come from list"
$rngK.Validation.ErrorMessage = "Value must be one of: protein."
